{"js": "// Replace the date line and the 25 two-digit-division answer cells with\n// their updated values, as described by the diff.\nconst replacements = [\n  [\"2023-11-13 Monday\", \"2023-11-14 Tuesday\"],\n  [\"34\u00f79=3, 7\", \"88\u00f78=11, 0\"],\n  [\"85\u00f76=14, 1\", \"40\u00f73=13, 1\"],\n  [\"39\u00f76=6, 3\", \"15\u00f75=3, 0\"],\n  [\"40\u00f74=10, 0\", \"10\u00f72=5, 0\"],\n  [\"33\u00f79=3, 6\", \"27\u00f79=3, 0\"],\n  [\"98\u00f77=14, 0\", \"49\u00f76=8, 1\"],\n  [\"67\u00f74=16, 3\", \"78\u00f78=9, 6\"],\n  [\"16\u00f74=4, 0\", \"18\u00f78=2, 2\"],\n  [\"88\u00f79=9, 7\", \"62\u00f72=31, 0\"],\n  [\"84\u00f77=12, 0\", \"75\u00f75=15, 0\"],\n  [\"65\u00f73=21, 2\", \"56\u00f74=14, 0\"],\n  [\"37\u00f74=9, 1\", \"33\u00f78=4, 1\"],\n  [\"15\u00f79=1, 6\", \"67\u00f75=13, 2\"],\n  [\"73\u00f72=36, 1\", \"11\u00f78=1, 3\"],\n  [\"53\u00f79=5, 8\", \"71\u00f79=7, 8\"],\n  [\"95\u00f79=10, 5\", \"28\u00f75=5, 3\"],\n  [\"42\u00f75=8, 2\", \"77\u00f75=15, 2\"],\n  [\"32\u00f76=5, 2\", \"11\u00f74=2, 3\"],\n  [\"50\u00f73=16, 2\", \"37\u00f78=4, 5\"],\n  [\"25\u00f72=12, 1\", \"79\u00f75=15, 4\"],\n  [\"10\u00f77=1, 3\", \"71\u00f74=17, 3\"],\n  [\"66\u00f78=8, 2\", \"30\u00f74=7, 2\"],\n  [\"10\u00f78=1, 2\", \"35\u00f79=3, 8\"],\n  [\"81\u00f79=9, 0\", \"27\u00f73=9, 0\"],\n  [\"17\u00f78=2, 1\", \"59\u00f77=8, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update master to output generated at aa3dc9e\n# Replace the date line and the 25 two-digit-division answer cells with\n# their updated values, as described by the diff.\n\n$replacements = @(\n    @(\"2023-11-13 Monday\", \"2023-11-14 Tuesday\"),\n    @(\"34\u00f79=3, 7\", \"88\u00f78=11, 0\"),\n    @(\"85\u00f76=14, 1\", \"40\u00f73=13, 1\"),\n    @(\"39\u00f76=6, 3\", \"15\u00f75=3, 0\"),\n    @(\"40\u00f74=10, 0\", \"10\u00f72=5, 0\"),\n    @(\"33\u00f79=3, 6\", \"27\u00f79=3, 0\"),\n    @(\"98\u00f77=14, 0\", \"49\u00f76=8, 1\"),\n    @(\"67\u00f74=16, 3\", \"78\u00f78=9, 6\"),\n    @(\"16\u00f74=4, 0\", \"18\u00f78=2, 2\"),\n    @(\"88\u00f79=9, 7\", \"62\u00f72=31, 0\"),\n    @(\"84\u00f77=12, 0\", \"75\u00f75=15, 0\"),\n    @(\"65\u00f73=21, 2\", \"56\u00f74=14, 0\"),\n    @(\"37\u00f74=9, 1\", \"33\u00f78=4, 1\"),\n    @(\"15\u00f79=1, 6\", \"67\u00f75=13, 2\"),\n    @(\"73\u00f72=36, 1\", \"11\u00f78=1, 3\"),\n    @(\"53\u00f79=5, 8\", \"71\u00f79=7, 8\"),\n    @(\"95\u00f79=10, 5\", \"28\u00f75=5, 3\"),\n    @(\"42\u00f75=8, 2\", \"77\u00f75=15, 2\"),\n    @(\"32\u00f76=5, 2\", \"11\u00f74=2, 3\"),\n    @(\"50\u00f73=16, 2\", \"37\u00f78=4, 5\"),\n    @(\"25\u00f72=12, 1\", \"79\u00f75=15, 4\"),\n    @(\"10\u00f77=1, 3\", \"71\u00f74=17, 3\"),\n    @(\"66\u00f78=8, 2\", \"30\u00f74=7, 2\"),\n    @(\"10\u00f78=1, 2\", \"35\u00f79=3, 8\"),\n    @(\"81\u00f79=9, 0\", \"27\u00f73=9, 0\"),\n    @(\"17\u00f78=2, 1\", \"59\u00f77=8, 3\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
